# Schulferien - unnötige Zeilen gelöscht
# Clear the leftover footnote / legend text block (rows 21-29) while
# keeping the existing cell formatting (number formats / alignment) intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21:G29").ClearContents()

# The explanatory picture/logo that used to sit above that block is no
# longer needed either.
if ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete()
}

# Update the view so the previously selected cell reflects where the user
# was working when the cleanup happened.
$ws.Range("A19:I33").Select() | Out-Null
